$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new column (T) with the 2023 data, matching the formatting --
# of the neighbouring 2022 column (S) / 2019 column (D) ------------------

# Row 2: thin divider cell above the new column - same look as Q2:S2.
$ws.Range("S2").Copy()
$ws.Range("T2").PasteSpecial(-4122)

# Row 3: year header "2023" - same look as the other year header cells
# (D3:P3), i.e. centred vertical alignment, bold, bottom border.
$ws.Range("D3").Copy()
$ws.Range("Q3:T3").PasteSpecial(-4122)
$ws.Range("Q3").Value = 2020
$ws.Range("R3").Value = 2021
$ws.Range("S3").Value = 2022
$ws.Range("T3").Value = 2023

# Row 4: absolute number of pensioners below subsistence level.
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Value = 263951

# Row 5: share of total population, percent.
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)
$ws.Range("T5").Value = 3.7

# Row 6: share of total pensioners, percent.
$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("T6").Value = 32.299999999999997

$ws.Application.CutCopyMode = $false

# --- Row 2 got a little taller in the refreshed layout ------------------
$ws.Rows(2).RowHeight = 16.5

# --- Restore the default selection (A1) instead of the stale C19 one ----
$ws.Range("A1").Select()
